$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1109388365702823

$ws.Range("B3").Value = -0.0016291604756874
$ws.Range("C3").Value = 0.0006523583492759668
$ws.Range("D3").Value = -0.174202558631421
$ws.Range("E3").Value = 0.05701432490148314
$ws.Range("F3").Value = -0.002907763688614172
$ws.Range("G3").Value = -0.0003505572627606288
$ws.Range("H3").Value = 0.1093096760945949

$ws.Range("B4").Value = 0.00139332588845911
$ws.Range("C4").Value = 0.0008484314525409661
$ws.Range("D4").Value = 3.987210621696983
$ws.Range("E4").Value = 0.008950513493478867
$ws.Range("F4").Value = -0.0002695747103406285
$ws.Range("G4").Value = 0.00305622648725885
$ws.Range("H4").Value = 0.1123321624587414

$ws.Range("B5").Value = 0.005581518874381719
$ws.Range("C5").Value = 0.001925553329321675
$ws.Range("D5").Value = 3.962111742747132
$ws.Range("E5").Value = 0.03831619871466767
$ws.Range("F5").Value = 0.001807490501437142
$ws.Range("G5").Value = 0.009355547247326297
$ws.Range("H5").Value = 0.116520355444664

$ws.Range("B6").Value = 0.01320294821679943
$ws.Range("C6").Value = 0.006764283919977927
$ws.Range("D6").Value = 6.753365811296738
$ws.Range("E6").Value = 1.239456779502655
$ws.Range("F6").Value = 0.002096001606625188
$ws.Range("G6").Value = 0.02081402781299595
$ws.Range("H6").Value = 0.1241417847870818

$ws.Range("B7").Value = 0.01354645026770222
$ws.Range("C7").Value = 0.004661625051623263
$ws.Range("D7").Value = 6.225105119279846
$ws.Range("E7").Value = 1.404637671028977
$ws.Range("F7").Value = 0.006553733490101049
$ws.Range("G7").Value = 0.01687766953044248
$ws.Range("H7").Value = 0.1244852868379845

$ws.Range("B8").Value = 0.01811551634441686
$ws.Range("C8").Value = 0.005306897619347763
$ws.Range("D8").Value = 5.833685441470382
$ws.Range("E8").Value = 1.510131500729628
$ws.Range("F8").Value = 0.01011036267972027
$ws.Range("G8").Value = 0.02197282198531864
$ws.Range("H8").Value = 0.1290543529146992

$ws.Range("B9").Value = 0.01973297011959893
$ws.Range("C9").Value = 0.005350861388470732
$ws.Range("D9").Value = 6.197859341183493
$ws.Range("E9").Value = 1.307030090128754
$ws.Range("F9").Value = 0.01098170142909386
$ws.Range("G9").Value = 0.0255692121485739
$ws.Range("H9").Value = 0.1306718066898812

$ws.Range("B10").Value = -0.1109388365702823
$ws.Range("C10").Value = 0.0004838750934690815
$ws.Range("D10").Value = -238.6626675955286
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.111887217531127
$ws.Range("G10").Value = -0.1099904556094376
$ws.Range("H10").Value = 0

$ws.Range("B11").Value = -0.05050962494755598
$ws.Range("C11").Value = 0.0005357100987420404
$ws.Range("D11").Value = -97.25040717689463
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -0.05155960098074986
$ws.Range("G11").Value = -0.04945964891436212
$ws.Range("H11").Value = 0.06042921162272634

$ws.Range("B12").Value = -0.04058149768070585
$ws.Range("C12").Value = 0.000512333044328639
$ws.Range("D12").Value = -81.71223266686891
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = -0.04158565537488565
$ws.Range("G12").Value = -0.03957733998652604
$ws.Range("H12").Value = 0.07035733888957647

$ws.Range("B13").Value = -0.03436092542434475
$ws.Range("C13").Value = 0.0005132955095480475
$ws.Range("D13").Value = -70.07177287373671
$ws.Range("E13").Value = [double]"1.143355490286354e-235"
$ws.Range("F13").Value = -0.03536696951093642
$ws.Range("G13").Value = -0.03335488133775304
$ws.Range("H13").Value = 0.07657791114593757

$ws.Range("B14").Value = -0.0253950015553659
$ws.Range("C14").Value = 0.0005082703755546098
$ws.Range("D14").Value = -52.50196569880657
$ws.Range("E14").Value = [double]"2.465802198087558e-68"
$ws.Range("F14").Value = -0.02639119651204766
$ws.Range("G14").Value = -0.02439880659868415
$ws.Range("H14").Value = 0.08554383501491641

$ws.Range("B15").Value = -0.02299253311204904
$ws.Range("C15").Value = 0.000502047578292078
$ws.Range("D15").Value = -46.34328035224348
$ws.Range("E15").Value = [double]"1.137325311142328e-21"
$ws.Range("F15").Value = -0.0239765315921851
$ws.Range("G15").Value = -0.022008534631913
$ws.Range("H15").Value = 0.08794630345823327

$ws.Range("B16").Value = -0.02140628938593895
$ws.Range("C16").Value = 0.0004935106785451833
$ws.Range("D16").Value = -43.97172109921745
$ws.Range("E16").Value = [double]"1.542454981320542e-08"
$ws.Range("F16").Value = -0.02237355578107833
$ws.Range("G16").Value = -0.02043902299079957
$ws.Range("H16").Value = 0.08953254718434338

$ws.Range("B17").Value = -0.01943525018266713
$ws.Range("C17").Value = 0.0004916961970893901
$ws.Range("D17").Value = -40.28367538507801
$ws.Range("E17").Value = [double]"2.488978164370008e-16"
$ws.Range("F17").Value = -0.0203989602432947
$ws.Range("G17").Value = -0.01847154012203956
$ws.Range("H17").Value = 0.09150358638761519

$ws.Range("B18").Value = -0.01685068517496338
$ws.Range("C18").Value = 0.0004950061612120125
$ws.Range("D18").Value = -35.19865291997507
$ws.Range("E18").Value = [double]"3.977962305757616e-08"
$ws.Range("F18").Value = -0.01782088266886523
$ws.Range("G18").Value = -0.01588048768106153
$ws.Range("H18").Value = 0.09408815139531894

$ws.Range("B19").Value = -0.01501808203926332
$ws.Range("C19").Value = 0.0004892297790235301
$ws.Range("D19").Value = -32.0160814389201
$ws.Range("E19").Value = 0.03984603746357507
$ws.Range("F19").Value = -0.01597695798478599
$ws.Range("G19").Value = -0.01405920609374064
$ws.Range("H19").Value = 0.095920754531019

$ws.Range("B20").Value = -0.01252108365273385
$ws.Range("C20").Value = 0.0004991999272394913
$ws.Range("D20").Value = -25.19818024462377
$ws.Range("E20").Value = 0.001746430920033519
$ws.Range("F20").Value = -0.01349950079166523
$ws.Range("G20").Value = -0.01154266651380248
$ws.Range("H20").Value = 0.09841775291754846

$ws.Range("B21").Value = -0.009419862408693969
$ws.Range("C21").Value = 0.0004994981427696024
$ws.Range("D21").Value = -19.51853334649007
$ws.Range("E21").Value = [double]"1.637029451479686e-05"
$ws.Range("F21").Value = -0.01039886401750905
$ws.Range("G21").Value = -0.008440860799878883
$ws.Range("H21").Value = 0.1015189741615884

$ws.Range("B22").Value = -0.006925118339445624
$ws.Range("C22").Value = 0.0004941934737623512
$ws.Range("D22").Value = -14.02053210846501
$ws.Range("E22").Value = 0.07372390267681711
$ws.Range("F22").Value = -0.007893722960439862
$ws.Range("G22").Value = -0.005956513718451388
$ws.Range("H22").Value = 0.1040137182308367

$ws.Range("B23").Value = -0.005915836906182538
$ws.Range("C23").Value = 0.0004867659838057682
$ws.Range("D23").Value = -11.56649070967306
$ws.Range("E23").Value = 0.02548451979662011
$ws.Range("F23").Value = -0.006869883862752247
$ws.Range("G23").Value = -0.004961789949612829
$ws.Range("H23").Value = 0.1050229996640998

$ws.Range("B24").Value = -0.005951234553062235
$ws.Range("C24").Value = 0.0004798332081502496
$ws.Range("D24").Value = -11.94808101774407
$ws.Range("E24").Value = 0.14750057447387
$ws.Range("F24").Value = -0.006891693499731146
$ws.Range("G24").Value = -0.005010775606393326
$ws.Range("H24").Value = 0.1049876020172201

$ws.Range("B25").Value = -0.002903862533398152
$ws.Range("C25").Value = 0.0004692235485583242
$ws.Range("D25").Value = -5.595519899335898
$ws.Range("E25").Value = 0.0903016206341075
$ws.Range("F25").Value = -0.003823526868278516
$ws.Range("G25").Value = -0.001984198198517787
$ws.Range("H25").Value = 0.1080349740368842

$ws.Range("B26").Value = 0.04076821288858381
$ws.Range("C26").Value = 0.005219080897656299
$ws.Range("D26").Value = 21.55276334085174
$ws.Range("E26").Value = 1.85880764553502
$ws.Range("F26").Value = 0.03405920785269608
$ws.Range("G26").Value = 0.04117867599054097
$ws.Range("H26").Value = 0.1517070494588661
